{"js": "const pairs = [\n  [\"11\u00d712=\", \"48\u00d725=\"],\n  [\"59\u00d798=\", \"20\u00d784=\"],\n  [\"56\u00d734=\", \"27\u00d733=\"],\n  [\"99\u00d785=\", \"17\u00d754=\"],\n  [\"92\u00d744=\", \"44\u00d739=\"],\n  [\"53\u00d778=\", \"81\u00d796=\"],\n  [\"65\u00d768=\", \"34\u00d735=\"],\n  [\"66\u00d786=\", \"85\u00d791=\"],\n  [\"17\u00d747=\", \"71\u00d758=\"],\n  [\"38\u00d721=\", \"87\u00d756=\"],\n  [\"38\u00d712=\", \"26\u00d750=\"],\n  [\"31\u00d764=\", \"68\u00d775=\"],\n  [\"30\u00d791=\", \"72\u00d724=\"],\n  [\"84\u00d732=\", \"27\u00d788=\"],\n  [\"28\u00d711=\", \"68\u00d713=\"],\n  [\"16\u00d738=\", \"40\u00d732=\"],\n  [\"18\u00d790=\", \"25\u00d721=\"],\n  [\"25\u00d716=\", \"43\u00d734=\"],\n  [\"85\u00d740=\", \"30\u00d724=\"],\n  [\"37\u00d711=\", \"85\u00d735=\"],\n  [\"34\u00d746=\", \"15\u00d711=\"],\n  [\"35\u00d747=\", \"94\u00d776=\"],\n  [\"42\u00d797=\", \"64\u00d736=\"],\n  [\"11\u00d734=\", \"93\u00d771=\"],\n  [\"62\u00d774=\", \"27\u00d794=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"11\u00d712=\", \"48\u00d725=\"),\n    @(\"59\u00d798=\", \"20\u00d784=\"),\n    @(\"56\u00d734=\", \"27\u00d733=\"),\n    @(\"99\u00d785=\", \"17\u00d754=\"),\n    @(\"92\u00d744=\", \"44\u00d739=\"),\n    @(\"53\u00d778=\", \"81\u00d796=\"),\n    @(\"65\u00d768=\", \"34\u00d735=\"),\n    @(\"66\u00d786=\", \"85\u00d791=\"),\n    @(\"17\u00d747=\", \"71\u00d758=\"),\n    @(\"38\u00d721=\", \"87\u00d756=\"),\n    @(\"38\u00d712=\", \"26\u00d750=\"),\n    @(\"31\u00d764=\", \"68\u00d775=\"),\n    @(\"30\u00d791=\", \"72\u00d724=\"),\n    @(\"84\u00d732=\", \"27\u00d788=\"),\n    @(\"28\u00d711=\", \"68\u00d713=\"),\n    @(\"16\u00d738=\", \"40\u00d732=\"),\n    @(\"18\u00d790=\", \"25\u00d721=\"),\n    @(\"25\u00d716=\", \"43\u00d734=\"),\n    @(\"85\u00d740=\", \"30\u00d724=\"),\n    @(\"37\u00d711=\", \"85\u00d735=\"),\n    @(\"34\u00d746=\", \"15\u00d711=\"),\n    @(\"35\u00d747=\", \"94\u00d776=\"),\n    @(\"42\u00d797=\", \"64\u00d736=\"),\n    @(\"11\u00d734=\", \"93\u00d771=\"),\n    @(\"62\u00d774=\", \"27\u00d794=\")\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2)\n}\n"}
